$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 237.11241462348252
$ws.Range("C2").Value = 173.28826462944508
$ws.Range("D2").Value = 238.09668928220884
$ws.Range("E2").Value = 175.68100970774961

$ws.Range("B3").Value = 217.36520251781573
$ws.Range("C3").Value = 171.82536395782364
$ws.Range("D3").Value = 211.33599447667382
$ws.Range("E3").Value = 177.64809922808044

$ws.Range("B1:E3").Select()
